$d = $word.ActiveDocument

# Locate, by content, the paragraph that reads
# "There is no documentation for these questions." and the empty paragraph
# that immediately precedes it (which sits right after the
# "Hi, If you have ended up here..." paragraph). Both paragraphs are to be
# removed in their entirety (including their own paragraph marks), so that
# the "Hi, If you have ended up here..." paragraph is directly followed by
# the remaining empty paragraph that used to come after
# "There is no documentation for these questions.".

$docParagraph = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*There is no documentation for these questions.*") {
        $docParagraph = $p
        break
    }
}

if ($docParagraph -ne $null) {
    # Delete the "There is no documentation..." paragraph, including its own
    # paragraph mark.
    $docRange = $d.Range($docParagraph.Range.Start, $docParagraph.Range.End)
    $docRange.Delete()

    # The paragraph immediately before it is the now-adjacent empty
    # paragraph; delete it too, including its own paragraph mark.
    $emptyParagraph = $docParagraph.Previous()
    $emptyRange = $d.Range($emptyParagraph.Range.Start, $emptyParagraph.Range.End)
    $emptyRange.Delete()
}
